$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "EmpresaTeste"
$ws.Range("B5").Value = "12.345.678/9123-45"
$ws.Range("C5").Value = "12345-678"
$ws.Range("D5").Value = "eg@gmail.com"
$ws.Range("E5").Value = "(71) 99945-2004"
$ws.Range("F5").Value = "Eg2024!"
